$wb = $excel.ActiveWorkbook

# The commit adds a new "login_code" column to the SCHUELER table.
$ws = $wb.Worksheets.Item("SCHUELER")

$ws.Range("D1").Value = "login_code"
$ws.Range("D2").Value = "aaa"
